$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: C15 "subNO" -> "questionCount"; D15 "TRQ_SCRSUB.SUBNO" -> cleared
$ws.Range("C15").Value = "questionCount"
$ws.Range("D15").Value = ""

# Row 16: clear all four cells (previously 答案描述 / 语音较低... / subDesc / TRQ_SCRSUB.SUBDESC)
$ws.Range("A16:D16").Value = ""

# Move the active selection from B16 to C21
$ws.Range("C21").Select() | Out-Null
